# docenteninfo.xlsx — replace long Dutch subject names in column C ("vakken")
# with short 2-letter subject codes, and fill the previously-empty
# "uren-onbeschikbaar"/"uren-liever-niet" columns (G/H) with the literal
# placeholder "none" for every teacher row that didn't already have a value
# there.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- columns G/H: fill previously-blank cells with "none" --------------
# (written first so the new shared string "none" lands at the same
# shared-strings index the original author's save produced)
$ws.Range("G2").Value  = "none"
$ws.Range("H2").Value  = "none"
$ws.Range("G4").Value  = "none"
$ws.Range("H4").Value  = "none"
$ws.Range("G5").Value  = "none"
$ws.Range("H5").Value  = "none"
$ws.Range("G6").Value  = "none"
$ws.Range("H6").Value  = "none"
$ws.Range("G7").Value  = "none"
$ws.Range("H7").Value  = "none"
$ws.Range("G8").Value  = "none"
$ws.Range("H8").Value  = "none"
$ws.Range("G9").Value  = "none"
$ws.Range("H9").Value  = "none"
$ws.Range("G10").Value = "none"
$ws.Range("H10").Value = "none"
$ws.Range("G11").Value = "none"
$ws.Range("H11").Value = "none"
$ws.Range("G12").Value = "none"
$ws.Range("H12").Value = "none"

# --- column C: vakken -> short codes -----------------------------------
$ws.Range("C2").Value  = "WI"   # wiskunde
$ws.Range("C3").Value  = "WI"   # wiskunde
$ws.Range("C4").Value  = "LO"   # gymnastiek
$ws.Range("C5").Value  = "HV"   # handenarbeid
$ws.Range("C6").Value  = "NL"   # nederlands
$ws.Range("C7").Value  = "LA|GR" # latijn|grieks
$ws.Range("C8").Value  = "GS"   # geschiedenis
$ws.Range("C9").Value  = "BI"   # biologie
$ws.Range("C10").Value = "AK"   # aardrijkskunde
$ws.Range("C11").Value = "EN"   # engels
$ws.Range("C12").Value = "ML"   # maatschappijleer

# --- cosmetic: page setup + selected cell, matching the saved view -----
$ws.PageSetup.PaperSize = 9   # xlPaperA4
$ws.PageSetup.Orientation = 1 # xlPortrait

$ws.Range("C13").Select() | Out-Null
